$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Update the underlying source values; dependent formulas (D19, E19, H19, I19)
# will recalculate automatically.
$ws.Range("B19").Value = 65320700
$ws.Range("F19").Value = 59209900

# Update the selected cell/active window on the sheet.
$ws.Range("F32").Select()
